$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$label = "Diferença 2024/02 - 2023/02"

# Row 2: Amapá
$ws.Range("A2").Value = "Amapá"
$ws.Range("B2").Value = $label
$ws.Range("C2").Value = 3.099832034109994
$ws.Range("D2").Value = "1º"

# Row 3: Pernambuco
$ws.Range("A3").Value = "Pernambuco"
$ws.Range("B3").Value = $label
$ws.Range("C3").Value = 2.651345123739418
$ws.Range("D3").Value = "2º"

# Row 4: Bahia
$ws.Range("A4").Value = "Bahia"
$ws.Range("B4").Value = $label
$ws.Range("C4").Value = 2.257958882804374
$ws.Range("D4").Value = "3º"

# Row 5: Piauí
$ws.Range("A5").Value = "Piauí"
$ws.Range("B5").Value = $label
$ws.Range("C5").Value = 2.172817929474661
$ws.Range("D5").Value = "4º"

# Row 6: Tocantins
$ws.Range("A6").Value = "Tocantins"
$ws.Range("B6").Value = $label
$ws.Range("C6").Value = 2.126230593957899
$ws.Range("D6").Value = "5º"

# Row 7: Nordeste
$ws.Range("A7").Value = "Nordeste"
$ws.Range("B7").Value = $label
$ws.Range("C7").Value = 1.898647905536663
$ws.Range("D7").Value = "6º"

# Row 8: Sergipe (now gains a D value)
$ws.Range("A8").Value = "Sergipe"
$ws.Range("B8").Value = $label
$ws.Range("C8").Value = 1.194889832701833
$ws.Range("D8").Value = "15º"

# Row 9: new row - Nordeste (duplicate, no D value)
$ws.Range("A9").Value = "Nordeste"
$ws.Range("B9").Value = $label
$ws.Range("C9").Value = 1.898647905536663

# Row 10: Brasil moves down one row (no D value)
$ws.Range("A10").Value = "Brasil"
$ws.Range("B10").Value = $label
$ws.Range("C10").Value = 1.143726150552311
